$wb = $excel.ActiveWorkbook

# Remember the originally active sheet so we can restore it at the end
# (selecting a range on another sheet would otherwise flip the active tab).
$origActive = $wb.ActiveSheet.Name

$ws = $wb.Worksheets.Item("Crypto")

# Rename "Crypto" -> "Built-in functions"
$ws.Name = "Built-in functions"

# New scenario rows appended to the sheet (rows 5-9), matching the style
# (font/fill/border) of the existing data rows (row 4).
$newRows = @(
    @{ Row = 5; Text = "Generate CityHash" },
    @{ Row = 6; Text = "Uncompress gzip data" },
    @{ Row = 7; Text = "Uncompress xz data" },
    @{ Row = 8; Text = "Generate random numbers" },
    @{ Row = 9; Text = "Return total system memory" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.Text
    $ws.Range("A4").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Row 5 also carries a (blank) B cell styled like the rest of column B.
$ws.Range("B5").Value = ""
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Update the sheet's remembered selection to C10, as in the authored edit.
[void]$ws.Range("C10").Select()

# Restore the original active sheet/tab.
[void]$wb.Worksheets.Item($origActive).Activate()

Write-Host "done"
